# Add team record columns (Wins, Losses, Ties) to the NYY 2006 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting (bold, centered, bordered) used by the
# existing header cells by copying the style from AC1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows (2-51): team record repeated for every player ---
$ws.Range("AD2:AD51").Value = 97
$ws.Range("AE2:AE51").Value = 65
$ws.Range("AF2:AF51").Value = 0
